$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7981066107749939
$ws.Range("B1").Value = 1.499279379844666
$ws.Range("C1").Value = 5.820702075958252
$ws.Range("D1").Value = 3.132373809814453
$ws.Range("E1").Value = 1.471651196479797
